$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -11.047
$ws.Range("D3").Value = -7.264999999999999
$ws.Range("A12").Value = -21.669
$ws.Range("C14").Value = -12.139
$ws.Range("C26").Value = -12.805
$ws.Range("D30").Value = -7.257
$ws.Range("C31").Value = -13.186
$ws.Range("A32").Value = -21.757
$ws.Range("C35").Value = -12.762
$ws.Range("A36").Value = -20.178
$ws.Range("C37").Value = -13.341
$ws.Range("A38").Value = -19.741
$ws.Range("D44").Value = -7.746
$ws.Range("C45").Value = -12.883
$ws.Range("A46").Value = -21.924
$ws.Range("A54").Value = -22.155
$ws.Range("A55").Value = -22.278
$ws.Range("C57").Value = -13.829
$ws.Range("D58").Value = -8.040000000000001
$ws.Range("A67").Value = -21.481
$ws.Range("A69").Value = -21.637
$ws.Range("A72").Value = -21.445
$ws.Range("D84").Value = -8.405999999999999
$ws.Range("D89").Value = -6.962000000000001
$ws.Range("A91").Value = -21.584
$ws.Range("D91").Value = -6.931999999999999
$ws.Range("A99").Value = -20.428
$ws.Range("C100").Value = -12.638
$ws.Range("C102").Value = -13.45
$ws.Range("D102").Value = -7.76
